# Data entry script: easily enter tools and employees, and record a
# tool checkout/checkin in the log.
#
# Adds a new employee (joe / 1234), a new tool (screw / 12345), logs a
# checkout+checkin of "screw" by joe, and records joe signing back in
# tool3 (previously checked out by emp3, never signed back in). Also
# clears the stray "Active" status left on tool3 and leaves the "tools"
# sheet as the active tab/selection, matching where the data entry
# finished.

$wb = $excel.ActiveWorkbook

$logSheet = $wb.Worksheets.Item("tool_checkout_log")
$employeesSheet = $wb.Worksheets.Item("employees")
$toolsSheet = $wb.Worksheets.Item("tools")

# --- tool_checkout_log: emp3's tool3 checkout gets signed back in by joe ---
$logSheet.Range("D4").Value = "02/10/2024 18:50"
$logSheet.Range("E4").Value = "joe"

# --- tool_checkout_log: joe checks out & immediately signs in "screw" ---
$logSheet.Range("A7").Value = "joe"
$logSheet.Range("B7").Value = "screw"
$logSheet.Range("C7").Value = "02/10/2024 18:49"
$logSheet.Range("D7").Value = "02/10/2024 18:49"
$logSheet.Range("E7").Value = "joe"

# --- employees: add new employee "joe" ---
$employeesSheet.Range("A12").Value = 1234
$employeesSheet.Range("B12").Value = "joe"

# --- tools: clear stray status, add new tools "hammer" and "screw" ---
$toolsSheet.Range("C4").ClearContents()
$toolsSheet.Range("A12").Value = 14543
$toolsSheet.Range("B12").Value = "hammer"
$toolsSheet.Range("A13").Value = 12345
$toolsSheet.Range("B13").Value = "screw"

# Leave the selection on the newly entered tool row, and make "tools"
# the active tab / sheet when the workbook is reopened.
$toolsSheet.Range("B12").Select()
$toolsSheet.Activate()
